# Update countries & provincias Spain
# Reorders several country rows in the "Pais" sheet (the sheet's shared-string
# table was reshuffled upstream), refreshes the "Datos actualizados" timestamp,
# and swaps the two pairs of numeric values whose rows changed position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer text update (04:50 -> 05:50)
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 29 de Marzo de 2020 a las 05:50'

# Country name (column A) reshuffles
$ws.Cells.Item(153, 1).Value = 'Dominica'
$ws.Cells.Item(154, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(155, 1).Value = 'Bahamas'
$ws.Cells.Item(156, 1).Value = 'Niger'
$ws.Cells.Item(159, 1).Value = 'Haiti'
$ws.Cells.Item(160, 1).Value = 'Laos'
$ws.Cells.Item(161, 1).Value = 'Mozambique'
$ws.Cells.Item(162, 1).Value = 'Guinea'
$ws.Cells.Item(164, 1).Value = 'Birmania'
$ws.Cells.Item(165, 1).Value = 'Surinam'
$ws.Cells.Item(172, 1).Value = 'Zimbabue'
$ws.Cells.Item(173, 1).Value = 'Gabon'
$ws.Cells.Item(174, 1).Value = 'Benin'
$ws.Cells.Item(176, 1).Value = 'Eritrea'
$ws.Cells.Item(178, 1).Value = 'San Bartolome'
$ws.Cells.Item(179, 1).Value = 'Siria'
$ws.Cells.Item(180, 1).Value = 'Angola'
$ws.Cells.Item(181, 1).Value = 'Montserrat'
$ws.Cells.Item(182, 1).Value = 'Fiyi'
$ws.Cells.Item(183, 1).Value = 'Mauritania'
$ws.Cells.Item(184, 1).Value = 'Sudan'
$ws.Cells.Item(185, 1).Value = 'Nepal'
$ws.Cells.Item(186, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(187, 1).Value = 'Congo'
$ws.Cells.Item(189, 1).Value = 'Libia'
$ws.Cells.Item(190, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(191, 1).Value = 'Republica del Chad'
$ws.Cells.Item(192, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(193, 1).Value = 'Liberia'
$ws.Cells.Item(194, 1).Value = 'Butan'
$ws.Cells.Item(199, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(200, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(201, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(203, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(204, 1).Value = 'Timor Oriental'

# Numeric values (columns D & H) that travelled with the rows that actually
# swapped places (Niger/Bahamas and Nepal/Sudan had different stats).
$ws.Cells.Item(155, 4).Value = 1
$ws.Cells.Item(155, 8).Value = 0
$ws.Cells.Item(156, 4).Value = 0
$ws.Cells.Item(156, 8).Value = 1

$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 8).Value = 1
$ws.Cells.Item(185, 4).Value = 1
$ws.Cells.Item(185, 8).Value = 0
